# Update numeric results on Sheet1 (case with 380 kV done)
# Columns B,C,D,F,G,H,I,J,K for rows 2-25 get refreshed result values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.4021169102478837
$ws.Cells.Item(2, 3).Value = 0.08249588546655673
$ws.Cells.Item(2, 4).Value = 0.2329709150594965
$ws.Cells.Item(2, 6).Value = 1.855214341059622
$ws.Cells.Item(2, 7).Value = 1.124136189207121
$ws.Cells.Item(2, 8).Value = 1.105830378752458
$ws.Cells.Item(2, 9).Value = 0.8963089761147387
$ws.Cells.Item(2, 10).Value = 0.3516043913743943
$ws.Cells.Item(2, 11).Value = 0.4669374073496328
$ws.Cells.Item(3, 2).Value = 0.3652612459559919
$ws.Cells.Item(3, 3).Value = 0.07445181431870651
$ws.Cells.Item(3, 4).Value = 0.2246172109262261
$ws.Cells.Item(3, 6).Value = 1.851231181785494
$ws.Cells.Item(3, 7).Value = 1.1232769715198
$ws.Cells.Item(3, 8).Value = 1.110521565025238
$ws.Cells.Item(3, 9).Value = 0.9010827201968681
$ws.Cells.Item(3, 10).Value = 0.3404906507548304
$ws.Cells.Item(3, 11).Value = 0.4234718100455268
$ws.Cells.Item(4, 2).Value = 0.3427365608944228
$ws.Cells.Item(4, 3).Value = 0.06953596095317494
$ws.Cells.Item(4, 4).Value = 0.2195927428210496
$ws.Cells.Item(4, 6).Value = 1.849812497413893
$ws.Cells.Item(4, 7).Value = 1.12344863650091
$ws.Cells.Item(4, 8).Value = 1.11390247276816
$ws.Cells.Item(4, 9).Value = 0.9044822053968886
$ws.Cells.Item(4, 10).Value = 0.3338768537891781
$ws.Cells.Item(4, 11).Value = 0.3969079353552445
$ws.Cells.Item(5, 2).Value = 0.3335843019469564
$ws.Cells.Item(5, 3).Value = 0.06753856571653216
$ws.Cells.Item(5, 4).Value = 0.2175716492593409
$ws.Cells.Item(5, 6).Value = 1.849492392089289
$ws.Cells.Item(5, 7).Value = 1.123694102252486
$ws.Cells.Item(5, 8).Value = 1.115406040713466
$ws.Cells.Item(5, 9).Value = 0.9059851726091921
$ws.Cells.Item(5, 10).Value = 0.3312344071027695
$ws.Cells.Item(5, 11).Value = 0.3861145151255414
$ws.Cells.Item(6, 2).Value = 0.3320662019248175
$ws.Cells.Item(6, 3).Value = 0.06720725374560743
$ws.Cells.Item(6, 4).Value = 0.2172376461841594
$ws.Cells.Item(6, 6).Value = 1.84945481608041
$ws.Cells.Item(6, 7).Value = 1.123745452853015
$ws.Cells.Item(6, 8).Value = 1.115663305296891
$ws.Cells.Item(6, 9).Value = 0.9062418410961222
$ws.Cells.Item(6, 10).Value = 0.3307988127238985
$ws.Cells.Item(6, 11).Value = 0.3843241933326738
$ws.Cells.Item(7, 2).Value = 0.3426130215948717
$ws.Cells.Item(7, 3).Value = 0.06950899966101076
$ws.Cells.Item(7, 4).Value = 0.2195653785725966
$ws.Cells.Item(7, 6).Value = 1.849807135952346
$ws.Cells.Item(7, 7).Value = 1.12345123672371
$ws.Cells.Item(7, 8).Value = 1.113922241018287
$ws.Cells.Item(7, 9).Value = 0.9045019987568175
$ws.Cells.Item(7, 10).Value = 0.3338410034581472
$ws.Cells.Item(7, 11).Value = 0.3967622432053304
$ws.Cells.Item(8, 2).Value = 0.3893875260666846
$ws.Cells.Item(8, 3).Value = 0.07971747048952693
$ws.Cells.Item(8, 4).Value = 0.230068863869576
$ws.Cells.Item(8, 6).Value = 1.853627653592042
$ws.Cells.Item(8, 7).Value = 1.123694612777442
$ws.Cells.Item(8, 8).Value = 1.107343998889903
$ws.Cells.Item(8, 9).Value = 0.8978576700586487
$ws.Cells.Item(8, 10).Value = 0.3477287069307522
$ws.Cells.Item(8, 11).Value = 0.4519248902967092
$ws.Cells.Item(9, 2).Value = 0.4819327897465655
$ws.Cells.Item(9, 3).Value = 0.0999212499359885
$ws.Cells.Item(9, 4).Value = 0.2514950329963597
$ws.Cells.Item(9, 6).Value = 1.86928147654973
$ws.Cells.Item(9, 7).Value = 1.129735770629594
$ws.Cells.Item(9, 8).Value = 1.098417549325504
$ws.Cells.Item(9, 9).Value = 0.8885509555989373
$ws.Cells.Item(9, 10).Value = 0.3766356797592465
$ws.Cells.Item(9, 11).Value = 0.5610745264093566
$ws.Cells.Item(10, 2).Value = 0.5504182571932006
$ws.Cells.Item(10, 3).Value = 0.1148802060982916
$ws.Cells.Item(10, 4).Value = 0.2677411478212832
$ws.Cells.Item(10, 6).Value = 1.885780324049605
$ws.Cells.Item(10, 7).Value = 1.137590338228179
$ws.Cells.Item(10, 8).Value = 1.094285658860514
$ws.Cells.Item(10, 9).Value = 0.8839923097492601
$ws.Cells.Item(10, 10).Value = 0.3989050972687949
$ws.Cells.Item(10, 11).Value = 0.6418576777275007
$ws.Cells.Item(11, 2).Value = 0.5816798665213696
$ws.Cells.Item(11, 3).Value = 0.1217110847289575
$ws.Cells.Item(11, 4).Value = 0.2752413607164499
$ws.Cells.Item(11, 6).Value = 1.894376445159565
$ws.Cells.Item(11, 7).Value = 1.141910732945959
$ws.Cells.Item(11, 8).Value = 1.092933778408678
$ws.Cells.Item(11, 9).Value = 0.8824154266537434
$ws.Cells.Item(11, 10).Value = 0.4092626362103857
$ws.Cells.Item(11, 11).Value = 0.6787361016440627
$ws.Cells.Item(12, 2).Value = 0.5935329906206732
$ws.Cells.Item(12, 3).Value = 0.1243015136452073
$ws.Cells.Item(12, 4).Value = 0.2780972341124084
$ws.Cells.Item(12, 6).Value = 1.897788752220208
$ws.Cells.Item(12, 7).Value = 1.143654609078226
$ws.Cells.Item(12, 8).Value = 1.092497811421225
$ws.Cells.Item(12, 9).Value = 0.8818899137893936
$ws.Cells.Item(12, 10).Value = 0.413217575230135
$ws.Cells.Item(12, 11).Value = 0.692719453421347
$ws.Cells.Item(13, 2).Value = 0.590979545886654
$ws.Cells.Item(13, 3).Value = 0.1237434530047494
$ws.Cells.Item(13, 4).Value = 0.2774814734306119
$ws.Cells.Item(13, 6).Value = 1.897046857100705
$ws.Cells.Item(13, 7).Value = 1.143274232141266
$ws.Cells.Item(13, 8).Value = 1.092588325164371
$ws.Cells.Item(13, 9).Value = 0.8819999044284685
$ws.Cells.Item(13, 10).Value = 0.4123643493859674
$ws.Cells.Item(13, 11).Value = 0.689707080626107
$ws.Cells.Item(14, 2).Value = 0.5826547311361026
$ws.Cells.Item(14, 3).Value = 0.121924126315804
$ws.Cells.Item(14, 4).Value = 0.2754760008315884
$ws.Cells.Item(14, 6).Value = 1.894654026625219
$ws.Cells.Item(14, 7).Value = 1.14205203934047
$ws.Cells.Item(14, 8).Value = 1.092896388354802
$ws.Cells.Item(14, 9).Value = 0.8823707560794745
$ws.Cells.Item(14, 10).Value = 0.4095873541466375
$ws.Cells.Item(14, 11).Value = 0.6798861558757494
$ws.Cells.Item(15, 2).Value = 0.5775574877781651
$ws.Cells.Item(15, 3).Value = 0.1208102207547768
$ws.Cells.Item(15, 4).Value = 0.2742496339306797
$ws.Cells.Item(15, 6).Value = 1.89320882237314
$ws.Cells.Item(15, 7).Value = 1.141317465543025
$ws.Cells.Item(15, 8).Value = 1.09309498046882
$ws.Cells.Item(15, 9).Value = 0.8826072453103819
$ws.Cells.Item(15, 10).Value = 0.4078906345173294
$ws.Cells.Item(15, 11).Value = 0.673872927375669
$ws.Cells.Item(16, 2).Value = 0.5483773619204442
$ws.Cells.Item(16, 3).Value = 0.1144343137687542
$ws.Cells.Item(16, 4).Value = 0.2672531933570497
$ws.Cells.Item(16, 6).Value = 1.885240520811593
$ws.Cells.Item(16, 7).Value = 1.137323060056758
$ws.Cells.Item(16, 8).Value = 1.094384633053949
$ws.Cells.Item(16, 9).Value = 0.8841053755035233
$ws.Cells.Item(16, 10).Value = 0.3982327864708708
$ws.Cells.Item(16, 11).Value = 0.6394501657024421
$ws.Cells.Item(17, 2).Value = 0.5305034858357658
$ws.Cells.Item(17, 3).Value = 0.1105295454550514
$ws.Cells.Item(17, 4).Value = 0.2629891623092533
$ws.Cells.Item(17, 6).Value = 1.88063179623056
$ws.Cells.Item(17, 7).Value = 1.135064292212761
$ws.Cells.Item(17, 8).Value = 1.095311010843005
$ws.Cells.Item(17, 9).Value = 0.8851518101418776
$ws.Cells.Item(17, 10).Value = 0.3923662465010693
$ws.Cells.Item(17, 11).Value = 0.6183658893984614
$ws.Cells.Item(18, 2).Value = 0.5202330206730608
$ws.Cells.Item(18, 3).Value = 0.1082860714794549
$ws.Cells.Item(18, 4).Value = 0.2605469415718744
$ws.Cells.Item(18, 6).Value = 1.878083624789241
$ws.Cells.Item(18, 7).Value = 1.133835426907879
$ws.Cells.Item(18, 8).Value = 1.095893507544744
$ws.Cells.Item(18, 9).Value = 0.8858004490424989
$ws.Cells.Item(18, 10).Value = 0.389013333999003
$ws.Cells.Item(18, 11).Value = 0.6062510068888116
$ws.Cells.Item(19, 2).Value = 0.5167573668945806
$ws.Cells.Item(19, 3).Value = 0.1075268901224149
$ws.Cells.Item(19, 4).Value = 0.2597218254531128
$ws.Cells.Item(19, 6).Value = 1.87723847801314
$ws.Cells.Item(19, 7).Value = 1.133431420072966
$ws.Cells.Item(19, 8).Value = 1.096099259368827
$ws.Cells.Item(19, 9).Value = 0.8860280924976252
$ws.Cells.Item(19, 10).Value = 0.3878817617006263
$ws.Cells.Item(19, 11).Value = 0.6021512285540211
$ws.Cells.Item(20, 2).Value = 0.5324051463537955
$ws.Cells.Item(20, 3).Value = 0.1109449616972142
$ws.Cells.Item(20, 4).Value = 0.2634420064202061
$ws.Cells.Item(20, 6).Value = 1.8811117772822
$ws.Cells.Item(20, 7).Value = 1.135297461523123
$ws.Cells.Item(20, 8).Value = 1.09520725529282
$ws.Cells.Item(20, 9).Value = 0.8850355749445953
$ws.Cells.Item(20, 10).Value = 0.3929885375789723
$ws.Cells.Item(20, 11).Value = 0.6206090829429058
$ws.Cells.Item(21, 2).Value = 0.5850995261564549
$ws.Cells.Item(21, 3).Value = 0.12245840570111
$ws.Cells.Item(21, 4).Value = 0.2760646311338348
$ws.Cells.Item(21, 6).Value = 1.89535259174238
$ws.Cells.Item(21, 7).Value = 1.142408097362733
$ws.Cells.Item(21, 8).Value = 1.092803840588644
$ws.Cells.Item(21, 9).Value = 0.8822598828999233
$ws.Cells.Item(21, 10).Value = 0.4104021348709779
$ws.Cells.Item(21, 11).Value = 0.682770305970763
$ws.Cells.Item(22, 2).Value = 0.6196257852151916
$ws.Cells.Item(22, 3).Value = 0.1300048051198814
$ws.Cells.Item(22, 4).Value = 0.2844057536617868
$ws.Cells.Item(22, 6).Value = 1.905575874166757
$ws.Cells.Item(22, 7).Value = 1.14768403580166
$ws.Cells.Item(22, 8).Value = 1.091675862290145
$ws.Cells.Item(22, 9).Value = 0.8808633343313801
$ws.Cells.Item(22, 10).Value = 0.4219739522318662
$ws.Cells.Item(22, 11).Value = 0.7235027173427966
$ws.Cells.Item(23, 2).Value = 0.6011905977173626
$ws.Cells.Item(23, 3).Value = 0.1259751679685337
$ws.Cells.Item(23, 4).Value = 0.2799455927567749
$ws.Cells.Item(23, 6).Value = 1.900035594670442
$ws.Cells.Item(23, 7).Value = 1.144810516033871
$ws.Cells.Item(23, 8).Value = 1.092237345724484
$ws.Cells.Item(23, 9).Value = 0.8815704387487173
$ws.Cells.Item(23, 10).Value = 0.4157803396744555
$ws.Cells.Item(23, 11).Value = 0.7017534418205287
$ws.Cells.Item(24, 2).Value = 0.5315453885850445
$ws.Cells.Item(24, 3).Value = 0.1107571475619693
$ws.Cells.Item(24, 4).Value = 0.2632372468509061
$ws.Cells.Item(24, 6).Value = 1.880894461865211
$ws.Cells.Item(24, 7).Value = 1.135191828493845
$ws.Cells.Item(24, 8).Value = 1.095254007708746
$ws.Cells.Item(24, 9).Value = 0.8850879783893788
$ws.Cells.Item(24, 10).Value = 0.3927071379976468
$ws.Cells.Item(24, 11).Value = 0.6195949140942787
$ws.Cells.Item(25, 2).Value = 0.4568098150937203
$ws.Cells.Item(25, 3).Value = 0.09443553596443621
$ws.Cells.Item(25, 4).Value = 0.2456100567082729
$ws.Cells.Item(25, 6).Value = 1.864170622281549
$ws.Cells.Item(25, 7).Value = 1.127503320700185
$ws.Cells.Item(25, 8).Value = 1.100406589609989
$ws.Cells.Item(25, 9).Value = 0.8906691391243342
$ws.Cells.Item(25, 10).Value = 0.3686352694350319
$ws.Cells.Item(25, 11).Value = 0.531442616872738
